$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Reservations sheet: refresh the reservation rows with newly suggested
# check-in / check-out dates (and the room/customer pairing that goes with
# them).
# ---------------------------------------------------------------------------
$res = $wb.Worksheets.Item("Reservations")

# Row 2: cong/dinh, room 101 (Single)
$res.Cells.Item(2,1).Value = "cong"
$res.Cells.Item(2,2).Value = "dinh"
$res.Cells.Item(2,3).Value = "cong@domain.com"
$res.Cells.Item(2,4).Value = "'101"
$res.Cells.Item(2,5).Value = 120.0
$res.Cells.Item(2,6).Value = "Single"
$res.Cells.Item(2,7).Value = "No"
$res.Cells.Item(2,8).Value = "01/01/2024 08:48:15"
$res.Cells.Item(2,9).Value = "01/04/2024 08:48:23"

# Row 3: cong/dinh, room 104 (Double)
$res.Cells.Item(3,1).Value = "cong"
$res.Cells.Item(3,2).Value = "dinh"
$res.Cells.Item(3,3).Value = "cong@domain.com"
$res.Cells.Item(3,4).Value = "'104"
$res.Cells.Item(3,5).Value = 300.0
$res.Cells.Item(3,6).Value = "Double"
$res.Cells.Item(3,7).Value = "No"
$res.Cells.Item(3,8).Value = "01/01/2024 08:47:15"
$res.Cells.Item(3,9).Value = "01/04/2024 08:47:20"

# Row 4: van/nguyen, room 104 (Double)
$res.Cells.Item(4,1).Value = "van"
$res.Cells.Item(4,2).Value = "nguyen"
$res.Cells.Item(4,3).Value = "van@domain.com"
$res.Cells.Item(4,4).Value = "'104"
$res.Cells.Item(4,5).Value = 300.0
$res.Cells.Item(4,6).Value = "Double"
$res.Cells.Item(4,7).Value = "No"
$res.Cells.Item(4,8).Value = "01/08/2024 09:00:35"
$res.Cells.Item(4,9).Value = "01/11/2024 09:00:41"

# ---------------------------------------------------------------------------
# Rooms sheet: the room list was trimmed back down to the two rooms that are
# actually referenced by a reservation (101, 104), dropping 102/103/106/107/
# 303/505 and reordering 104 ahead of 101.
# ---------------------------------------------------------------------------
$rooms = $wb.Worksheets.Item("Rooms")

# Delete rows 4 through 9 (rooms 103, 106, 107, 303, 505), bottom-up so the
# row indices of the rows still to be removed don't shift underneath us.
for ($r = 9; $r -ge 4; $r--) {
    $rooms.Rows.Item($r).EntireRow.Delete()
}

# Re-populate the remaining two data rows in the new order: 104 then 101.
$rooms.Cells.Item(2,1).Value = "'104"
$rooms.Cells.Item(2,2).Value = 300.0
$rooms.Cells.Item(2,3).Value = "Double"
$rooms.Cells.Item(2,4).Value = "No"

$rooms.Cells.Item(3,1).Value = "'101"
$rooms.Cells.Item(3,2).Value = 120.0
$rooms.Cells.Item(3,3).Value = "Single"
$rooms.Cells.Item(3,4).Value = "No"
